$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the serverTarget values in column A (rows 2-6) from "test17" to "test14"
$ws.Range("A2").Value = "test14"
$ws.Range("A3").Value = "test14"
$ws.Range("A4").Value = "test14"
$ws.Range("A5").Value = "test14"
$ws.Range("A6").Value = "test14"

# Move the active selection to A2 (was B10)
$ws.Range("A2").Select()
